$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the now-unused sample script names, keep DemoTest.xlsx in B2
$ws.Range("B3").Value = $null
$ws.Range("B4").Value = $null
$ws.Range("B5").Value = $null

# Move the active selection to B4
$ws.Range("B4").Select()
